# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3434
$ws1.Range("F5").Value = 1696
$ws1.Range("F6").Value = 86
$ws1.Range("F7").Value = 334

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3434
$ws4.Range("F5").Value = 1696
$ws4.Range("F6").Value = 86
$ws4.Range("F8").Value = 334
